# Generate Report for Handback
#
# A second handback file (c227cf0a-e2a6-4b58-96a5-4d13971bb925) has landed
# alongside the existing one (3c082223-2ec7-449d-b886-679c07488fc5, itself
# the replacement for the former d93ad8b6-0682-451b-9c8f-ad09f23c6808 run).
# Refresh the existing report row with the latest run's data and append a
# new row for the newly generated file on every sheet (Overview, zh-cn,
# de-de), keeping each sheet's Excel table in sync.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop the stale hyperlinks so the display text + targets can be rebuilt
# from scratch (an in-place update leaves the old display text behind).
$wsOverview.Hyperlinks.Delete()

# Row 2: update in place to point at the latest handback run.
$wsOverview.Range("A2").Value = "3c082223-2ec7-449d-b886-679c07488fc5.md"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G2").Value = "2016-08-15 18:55:52"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54892cd24efde05d1459be5a8ddb69ea1edcfc42/e2e/3c082223-2ec7-449d-b886-679c07488fc5.md", [Type]::Missing, [Type]::Missing, "e2e\3c082223-2ec7-449d-b886-679c07488fc5.md") | Out-Null

# Row 3 (new): the second file handed back in the same run.
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "c227cf0a-e2a6-4b58-96a5-4d13971bb925.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G3").Value = "2016-08-15 18:55:52"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54892cd24efde05d1459be5a8ddb69ea1edcfc42/e2e/c227cf0a-e2a6-4b58-96a5-4d13971bb925.md", [Type]::Missing, [Type]::Missing, "e2e\c227cf0a-e2a6-4b58-96a5-4d13971bb925.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Hyperlinks.Delete()

# Row 2: update in place.
$wsZhCn.Range("A2").Value = "3c082223-2ec7-449d-b886-679c07488fc5.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D2").Value = "e2e"
$wsZhCn.Range("E2").Value = "ht"
$wsZhCn.Range("F2").Value = "False"
$wsZhCn.Range("G2").Value = "3c082223-2ec7-449d-b886-679c07488fc5.590132d29df7c28a29422d4ccefc71c2df71c1dc.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-15 18:55:47"
$wsZhCn.Range("J2").Value = "3c082223-2ec7-449d-b886-679c07488fc5.590132d29df7c28a29422d4ccefc71c2df71c1dc.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-15 18:56:09"
$wsZhCn.Range("L2").Value = ""
$wsZhCn.Range("M2").Value = "True"
$wsZhCn.Range("N2").Value = ""
$wsZhCn.Range("O2").Value = "False"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54892cd24efde05d1459be5a8ddb69ea1edcfc42/e2e/3c082223-2ec7-449d-b886-679c07488fc5.md", [Type]::Missing, [Type]::Missing, "3c082223-2ec7-449d-b886-679c07488fc5.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54892cd24efde05d1459be5a8ddb69ea1edcfc42/e2e/3c082223-2ec7-449d-b886-679c07488fc5.md", [Type]::Missing, [Type]::Missing, "3c082223-2ec7-449d-b886-679c07488fc5.md") | Out-Null

# Row 3 (new): the second file handed back in the same run.
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "c227cf0a-e2a6-4b58-96a5-4d13971bb925.f32646b5f3b3172a53687b08c5c59270abf9c234.zh-cn.xlf"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").Value = "2016-08-15 18:55:47"
$wsZhCn.Range("J3").Value = "c227cf0a-e2a6-4b58-96a5-4d13971bb925.f32646b5f3b3172a53687b08c5c59270abf9c234.zh-cn.xlf"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").Value = "2016-08-15 18:56:09"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54892cd24efde05d1459be5a8ddb69ea1edcfc42/e2e/c227cf0a-e2a6-4b58-96a5-4d13971bb925.md", [Type]::Missing, [Type]::Missing, "c227cf0a-e2a6-4b58-96a5-4d13971bb925.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54892cd24efde05d1459be5a8ddb69ea1edcfc42/e2e/c227cf0a-e2a6-4b58-96a5-4d13971bb925.md", [Type]::Missing, [Type]::Missing, "c227cf0a-e2a6-4b58-96a5-4d13971bb925.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Hyperlinks.Delete()

# Row 2: update in place.
$wsDeDe.Range("A2").Value = "3c082223-2ec7-449d-b886-679c07488fc5.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D2").Value = "e2e"
$wsDeDe.Range("E2").Value = "ht"
$wsDeDe.Range("F2").Value = "False"
$wsDeDe.Range("G2").Value = "3c082223-2ec7-449d-b886-679c07488fc5.590132d29df7c28a29422d4ccefc71c2df71c1dc.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-15 18:55:52"
$wsDeDe.Range("J2").Value = "3c082223-2ec7-449d-b886-679c07488fc5.590132d29df7c28a29422d4ccefc71c2df71c1dc.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-15 18:56:17"
$wsDeDe.Range("L2").Value = ""
$wsDeDe.Range("M2").Value = "True"
$wsDeDe.Range("N2").Value = ""
$wsDeDe.Range("O2").Value = "False"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54892cd24efde05d1459be5a8ddb69ea1edcfc42/e2e/3c082223-2ec7-449d-b886-679c07488fc5.md", [Type]::Missing, [Type]::Missing, "3c082223-2ec7-449d-b886-679c07488fc5.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54892cd24efde05d1459be5a8ddb69ea1edcfc42/e2e/3c082223-2ec7-449d-b886-679c07488fc5.md", [Type]::Missing, [Type]::Missing, "3c082223-2ec7-449d-b886-679c07488fc5.md") | Out-Null

# Row 3 (new): the second file handed back in the same run.
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "c227cf0a-e2a6-4b58-96a5-4d13971bb925.f32646b5f3b3172a53687b08c5c59270abf9c234.de-de.xlf"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").Value = "2016-08-15 18:55:52"
$wsDeDe.Range("J3").Value = "c227cf0a-e2a6-4b58-96a5-4d13971bb925.f32646b5f3b3172a53687b08c5c59270abf9c234.de-de.xlf"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").Value = "2016-08-15 18:56:17"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54892cd24efde05d1459be5a8ddb69ea1edcfc42/e2e/c227cf0a-e2a6-4b58-96a5-4d13971bb925.md", [Type]::Missing, [Type]::Missing, "c227cf0a-e2a6-4b58-96a5-4d13971bb925.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54892cd24efde05d1459be5a8ddb69ea1edcfc42/e2e/c227cf0a-e2a6-4b58-96a5-4d13971bb925.md", [Type]::Missing, [Type]::Missing, "c227cf0a-e2a6-4b58-96a5-4d13971bb925.md") | Out-Null

Write-Host "Handback status report regenerated: added c227cf0a-e2a6-4b58-96a5-4d13971bb925 alongside 3c082223-2ec7-449d-b886-679c07488fc5 on Overview, zh-cn and de-de."
